# Applies scheduled-runner price/profit updates to Hyperion_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 1058
$ws.Range("I41").Value = 869.6667
$ws.Range("K41").Value = 869.6667
$ws.Range("M41").Value = -429.6667

# Row 43
$ws.Range("H43").Value = 1828.2858
$ws.Range("I43").Value = 1851
$ws.Range("J43").Value = 1692
$ws.Range("K43").Value = 1851
$ws.Range("L43").Value = 1692
$ws.Range("M43").Value = -1782
$ws.Range("N43").Value = -1830

# Row 70
$ws.Range("H70").Value = 6947.372
$ws.Range("I70").Value = 3199.5
$ws.Range("K70").Value = 9598.5
$ws.Range("M70").Value = -9328.5

# Row 73
$ws.Range("H73").Value = 6947.372
$ws.Range("I73").Value = 3199.5
$ws.Range("K73").Value = 9598.5
$ws.Range("M73").Value = -8662.5

# Row 74
$ws.Range("H74").Value = 7411.591
$ws.Range("J74").Value = 7742.9
$ws.Range("L74").Value = 7742.9
$ws.Range("N74").Value = -9614.9

# Row 77
$ws.Range("H77").Value = 7411.591
$ws.Range("J77").Value = 7742.9
$ws.Range("L77").Value = 38714.5
$ws.Range("N77").Value = -48074.5

# Row 80
$ws.Range("H80").Value = 4931.706
$ws.Range("I80").Value = 1333.1666
$ws.Range("K80").Value = 3999.4998
$ws.Range("M80").Value = -3001.4998

# Row 83
$ws.Range("H83").Value = 4931.706
$ws.Range("I83").Value = 1333.1666
$ws.Range("K83").Value = 11998.4994
$ws.Range("M83").Value = -7006.499400000001

# Row 98
$ws.Range("H98").Value = 1416.7097
$ws.Range("I98").Value = 1263.9333
$ws.Range("K98").Value = 1263.9333
$ws.Range("M98").Value = 234.0667000000001

# Row 122
$ws.Range("H122").Value = 1416.7097
$ws.Range("I122").Value = 1263.9333
$ws.Range("K122").Value = 3791.7999
$ws.Range("M122").Value = -1341.7999

# Row 125
$ws.Range("H125").Value = 3108.4807
$ws.Range("J125").Value = 3522.7
$ws.Range("L125").Value = 31704.3
$ws.Range("N125").Value = -36624.3

# Row 132
$ws.Range("H132").Value = 2081
$ws.Range("I132").Value = 2083.7932
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 6251.3796
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -3721.3796
$ws.Range("N132").Value = -11060

# Row 135
$ws.Range("H135").Value = 1013.4762
$ws.Range("I135").Value = 751.94116
$ws.Range("J135").Value = 2125
$ws.Range("K135").Value = 6767.47044
$ws.Range("L135").Value = 19125
$ws.Range("M135").Value = -4232.47044
$ws.Range("N135").Value = -24195

# Row 137
$ws.Range("H137").Value = 2590.8
$ws.Range("I137").Value = 2291.4583
$ws.Range("J137").Value = 2932.9048
$ws.Range("K137").Value = 6874.374899999999
$ws.Range("L137").Value = 8798.714399999999
$ws.Range("M137").Value = -4324.374899999999
$ws.Range("N137").Value = -13898.7144

# Row 138
$ws.Range("H138").Value = 2624.17
$ws.Range("J138").Value = 2999.5732
$ws.Range("L138").Value = 8998.7196
$ws.Range("N138").Value = -19278.7196


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4051.8738
$ws.Range("I32").Value = 2657.4343
$ws.Range("K32").Value = 2657.4343
$ws.Range("M32").Value = -2370.4343

# Row 45
$ws.Range("H45").Value = 67130.56
$ws.Range("I45").Value = 112925.89
$ws.Range("K45").Value = 112925.89
$ws.Range("M45").Value = -112548.89

# Row 61
$ws.Range("H61").Value = 1746.12
$ws.Range("I61").Value = 1658.8695
$ws.Range("J61").Value = 2749.5
$ws.Range("K61").Value = 1658.8695
$ws.Range("L61").Value = 2749.5
$ws.Range("M61").Value = -1446.8695
$ws.Range("N61").Value = -3173.5

# Row 88
$ws.Range("H88").Value = 1370.3334
$ws.Range("J88").Value = 1800
$ws.Range("L88").Value = 1800
$ws.Range("N88").Value = -2612

# Row 91
$ws.Range("H91").Value = 1370.3334
$ws.Range("J91").Value = 1800
$ws.Range("L91").Value = 1800
$ws.Range("N91").Value = -4608

# Row 102
$ws.Range("H102").Value = 5199.5713
$ws.Range("I102").Value = 4774.875
$ws.Range("K102").Value = 4774.875
$ws.Range("M102").Value = -3152.875

# Row 132
$ws.Range("H132").Value = 2945.2632
$ws.Range("I132").Value = 2247.8125
$ws.Range("J132").Value = 6665
$ws.Range("K132").Value = 6743.4375
$ws.Range("L132").Value = 19995
$ws.Range("M132").Value = -4213.4375
$ws.Range("N132").Value = -25055

# Row 136
$ws.Range("H136").Value = 1746.12
$ws.Range("I136").Value = 1658.8695
$ws.Range("J136").Value = 2749.5
$ws.Range("K136").Value = 4976.6085
$ws.Range("L136").Value = 8248.5
$ws.Range("M136").Value = -2426.6085
$ws.Range("N136").Value = -13348.5


$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 11211.218
$ws.Range("I86").Value = 8795.166999999999
$ws.Range("J86").Value = 19909
$ws.Range("K86").Value = 8795.166999999999
$ws.Range("L86").Value = 19909
$ws.Range("M86").Value = -7672.166999999999
$ws.Range("N86").Value = -22155

# Row 89
$ws.Range("H89").Value = 11211.218
$ws.Range("I89").Value = 8795.166999999999
$ws.Range("J89").Value = 19909
$ws.Range("K89").Value = 43975.835
$ws.Range("L89").Value = 99545
$ws.Range("M89").Value = -38359.835
$ws.Range("N89").Value = -110777

# Row 94
$ws.Range("H94").Value = 4768.372
$ws.Range("I94").Value = 699.64514
$ws.Range("J94").Value = 15279.25
$ws.Range("K94").Value = 699.64514
$ws.Range("L94").Value = 15279.25
$ws.Range("M94").Value = -248.64514
$ws.Range("N94").Value = -16181.25

# Row 134
$ws.Range("H134").Value = 2550.2954
$ws.Range("I134").Value = 732.0323
$ws.Range("K134").Value = 2196.0969
$ws.Range("M134").Value = 338.9031


$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1719.2307
$ws.Range("I58").Value = 1357
$ws.Range("K58").Value = 1357
$ws.Range("M58").Value = -1154

# Row 134
$ws.Range("H134").Value = 30666.666
$ws.Range("I134").Value = 45595.047
$ws.Range("K134").Value = 136785.141
$ws.Range("M134").Value = -134250.141

# Row 136
$ws.Range("H136").Value = 1719.2307
$ws.Range("I136").Value = 1357
$ws.Range("K136").Value = 4071
$ws.Range("M136").Value = -1521


$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 6683918.5
$ws.Range("J4").Value = 2434218.2
$ws.Range("L4").Value = 7302654.600000001
$ws.Range("N4").Value = -7302878.600000001

# Row 12
$ws.Range("H12").Value = 80
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 88.77778000000001
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 266.33334
$ws.Range("M12").Value = 170
$ws.Range("N12").Value = -612.33334

# Row 40
$ws.Range("H40").Value = 19.461538
$ws.Range("I40").Value = 20.11111
$ws.Range("J40").Value = 18
$ws.Range("K40").Value = 80.44444
$ws.Range("L40").Value = 72
$ws.Range("M40").Value = -11.44444
$ws.Range("N40").Value = -210

# Row 109
$ws.Range("H109").Value = 669.2
$ws.Range("I109").Value = 669.2
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 2007.6
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -967.6000000000001
$ws.Range("N109").ClearContents()

# Row 137
$ws.Range("H137").Value = 3505.7646
$ws.Range("I137").Value = 1745.8
$ws.Range("J137").Value = 4239.0835
$ws.Range("K137").Value = 5237.4
$ws.Range("L137").Value = 12717.2505
$ws.Range("M137").Value = -137.3999999999996
$ws.Range("N137").Value = -22917.2505

# Row 138
$ws.Range("H138").Value = 3266.875
$ws.Range("I138").Value = 2876.4285
$ws.Range("K138").Value = 8629.2855
$ws.Range("M138").Value = -3489.2855


$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 70258.07000000001
$ws.Range("I102").Value = 2573.8572
$ws.Range("J102").Value = 129481.75
$ws.Range("K102").Value = 2573.8572
$ws.Range("L102").Value = 129481.75
$ws.Range("M102").Value = -951.8571999999999
$ws.Range("N102").Value = -132725.75

# Row 122
$ws.Range("H122").Value = 82727.61
$ws.Range("I122").Value = 125149.07
$ws.Range("J122").Value = 3187.375
$ws.Range("K122").Value = 375447.21
$ws.Range("L122").Value = 9562.125
$ws.Range("M122").Value = -372997.21
$ws.Range("N122").Value = -14462.125


$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1161.7693
$ws.Range("I16").Value = 898
$ws.Range("J16").Value = 2269.6
$ws.Range("K16").Value = 898
$ws.Range("L16").Value = 2269.6
$ws.Range("M16").Value = -728
$ws.Range("N16").Value = -2609.6

# Row 40
$ws.Range("H40").Value = 7364.769
$ws.Range("I40").Value = 4142.2856
$ws.Range("K40").Value = 4142.2856
$ws.Range("M40").Value = -4006.2856

# Row 55
$ws.Range("H55").Value = 1993.3334
$ws.Range("I55").Value = 2456.25
$ws.Range("J55").Value = 1464.2858
$ws.Range("K55").Value = 2456.25
$ws.Range("L55").Value = 1464.2858
$ws.Range("M55").Value = -2283.25
$ws.Range("N55").Value = -1810.2858

# Row 93
$ws.Range("H93").Value = 2328.6924
$ws.Range("I93").Value = 2414.3333
$ws.Range("J93").Value = 1301
$ws.Range("K93").Value = 2414.3333
$ws.Range("L93").Value = 1301
$ws.Range("M93").Value = -1166.3333
$ws.Range("N93").Value = -3797


$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 77513.07000000001
$ws.Range("J41").Value = 77513.07000000001
$ws.Range("L41").Value = 77513.07000000001
$ws.Range("N41").Value = -78293.07000000001

# Row 132
$ws.Range("H132").Value = 341051.34
$ws.Range("I132").Value = 8893.096
$ws.Range("J132").Value = 1116087.2
$ws.Range("K132").Value = 26679.288
$ws.Range("L132").Value = 3348261.6
$ws.Range("M132").Value = -24149.288
$ws.Range("N132").Value = -3353321.6
